$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.199.62"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "1.855.96"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'0.7061"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'237.52"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.07948"
$ws.Range("E8").Value = "  +4.14%  "
$ws.Range("D9").Value = "'0.3017"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'23.33"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'0.08173"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.862.93"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'5.155"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "'0.7004"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "'89.59"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "29.175.38"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'5.792"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "'0.000007840"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'13.21"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'235.59"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "2.087.79"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'7.346"
$ws.Range("E24").Value = "  -3.65%  "
$ws.Range("D25").Value = "'161.79"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'8.894"
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "'0.1424"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").Value = "'18.00"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'1.921"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").Value = "'1.429"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Value = "'4.353"
$ws.Range("E32").Value = "  -3.39%  "
$ws.Range("D33").Value = "'4.012"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'0.05191"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'1.159"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").Value = "'0.6984"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -2.84%  "
$ws.Range("D38").Value = "'2.672"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'0.01837"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "'0.9295"
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("D42").Value = "1.121.06"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("D43").Value = "'0.4238"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "'5.829"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("D45").Value = "'69.43"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'102.59"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").Value = "'1.753"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "1.992.95"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'9.083"
$ws.Range("D51").Value = "'0.05925"
$ws.Range("E51").Value = "  +0.91%  "
